$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.625.54"
$ws.Range("D3").Value = "3.447.12"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.03%  "
$ws.Range("D7").Value = "3.447.62"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "4.038.40"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.15%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "3.450.97"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "61.736.39"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("E19").Value = "  +8.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").Value = "3.596.86"
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000122"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  +7.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.49%  "
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.796"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "2.593.64"
$ws.Range("E49").Value = "  +10.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.62%  "
$ws.Range("E51").Value = "  +0.57%  "
